$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows extracted from the scraped product listing, but this time also
# keeping entries whose parsed price fell at or below R$10,00 (those show up
# with the raw, unparsed "const integers = ..." debug text instead of a
# clean price — a bug in the extractor that slipped past the >10 filter).
$data = @(
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL", "3.699.00"),
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL", ".00"),
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL", "3.699.00"),
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL", "3.699.00"),
    @("91697550", "Ar Condicionado Split 24000 BTUs Quente e Frio Branco 220V Series A1 TCL", "const integers = '3.699'.00"),
    @("90163990", "Ventilador de Teto com Controle de Parede Fenix 3 Pás 96 cm 127V (110V) Ventisol", "const integers = '219'.90"),
    @("90163990", "Ventilador de Teto com Controle de Parede Fenix 3 Pás 96 cm 127V (110V) Ventisol", "const integers = '219'.90"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.9998.999"),
    @("91989296", "Banheira de Imersão Zen 150x72cm Branco Sensea", "8.999.00"),
    @("89062981", "Tapete de Banheiro Antiderrapante Cocoon Poliéster Grafite 40x60cm Sensea", "99.90"),
    @("1571352810", "Lâmpada Filamento E14 Ba35 Velachama 2w 127v - Foxlux", "11.90"),
    @("9043764190437655", "Ar Condicionado Cassete Atualle Eco Frio 60000BTUs 220V Elgin", "13.049.80"),
    @("1571352810", "Lâmpada Filamento E14 Ba35 Velachama 2w 127v - Foxlux", "11.90")
)

$startRow = 54
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]

    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $data[$i][1]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $data[$i][2]
}
